# ProcessingTimeExtrapolation.xlsx - add a new "Sheet1" worksheet containing
# a Species ID x Set Size probability lookup table, inserted right before the
# "Forest1 Processing Automated" sheet.

$wb = $excel.ActiveWorkbook

$target = $wb.Worksheets.Item("Forest1 Processing Automated")
$newSheet = $wb.Worksheets.Add($target)
$newSheet.Name = "Sheet1"

# Header row: "Set Size" label over the set-size columns
$newSheet.Range("B2").Value = "Set Size"

# Column headers for set sizes 3-6
$newSheet.Range("A3").Value = "Species ID"
$newSheet.Range("B3").Value = 3
$newSheet.Range("C3").Value = 4
$newSheet.Range("D3").Value = 5
$newSheet.Range("E3").Value = 6

# Species ID rows 1-9
$newSheet.Range("A4").Value = 1
$newSheet.Range("A5").Value = 2
$newSheet.Range("A6").Value = 3
$newSheet.Range("A7").Value = 4
$newSheet.Range("A8").Value = 5
$newSheet.Range("A9").Value = 6
$newSheet.Range("A10").Value = 7
$newSheet.Range("A11").Value = 8
$newSheet.Range("A12").Value = 9

# Known probability values for species 1 and 2
$newSheet.Range("B4").Value = 0.2
$newSheet.Range("C4").Value = 0.3
$newSheet.Range("B5").Value = 0.17
$newSheet.Range("C5").Value = 0.05

# Formatting: center align the whole table, bold the header row/column
$newSheet.Range("A2:E12").HorizontalAlignment = -4108
$newSheet.Range("A3:E3").Font.Bold = $true
$newSheet.Range("A4:A12").Font.Bold = $true

# Column width for Species ID column
$newSheet.Columns.Item(1).ColumnWidth = 10

# Selection on the new sheet
$newSheet.Range("C6").Select()

# Restore selection/active-cell bookkeeping on "Forest1 Processing Automated"
$target.Range("J31").Select()
